# Updated contribution sheet by parthey for week 2 iteration
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contributions")

# Iteration 2 contribution grades & comments (rows 14-18)
$ws.Range("C14").Value = "HD"
$ws.Range("D14").Value = "Good work on given user story tasks producing tests and functionality."

$ws.Range("C15").Value = "HD"
$ws.Range("D15").Value = "Keeping the team motivated, assigning work to individual members. Good work on the admin page."

$ws.Range("C16").Value = "HD"
$ws.Range("D16").Value = "Good work on creating adding and deleting functionality for parties and candidates."

$ws.Range("C17").Value = "HD"
$ws.Range("D17").Value = "Amazing work on the frontend side of voting page. "

$ws.Range("C18").Value = "P"
$ws.Range("D18").Value = "Some commitment to work in a group but no visible work done yet."

# Row 15 grows to fit the wrapped comment text
$ws.Rows.Item(15).RowHeight = 18.75

# Leave the cursor where the author left off
$ws.Range("D22").Select() | Out-Null
